$wb = $excel.ActiveWorkbook

# This workbook is a Final Fantasy XIV Leve-profit tracker (Brynhildr server).
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) is a crafting-job
# table of market-board prices/profits (columns H:N) that gets refreshed by a
# scheduled data-pull. Apply the refreshed values cell by cell.


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 819.34326
$ws.Range("I17").Value = 1063
$ws.Range("J17").Value = 781.5345
$ws.Range("K17").Value = 3189
$ws.Range("L17").Value = 2344.6035
$ws.Range("M17").Value = -3021
$ws.Range("N17").Value = -2680.6035
$ws.Range("H28").Value = 769.55554
$ws.Range("I28").Value = 816.1667
$ws.Range("J28").Value = 676.3333
$ws.Range("K28").Value = 816.1667
$ws.Range("L28").Value = 676.3333
$ws.Range("M28").Value = -331.1667
$ws.Range("N28").Value = -1646.3333
$ws.Range("H86").Value = 4947.875
$ws.Range("I86").Value = 7124.4
$ws.Range("J86").Value = 1320.3334
$ws.Range("K86").Value = 7124.4
$ws.Range("L86").Value = 1320.3334
$ws.Range("M86").Value = -6001.4
$ws.Range("N86").Value = -3566.3334
$ws.Range("H89").Value = 4947.875
$ws.Range("I89").Value = 7124.4
$ws.Range("J89").Value = 1320.3334
$ws.Range("K89").Value = 35622
$ws.Range("L89").Value = 6601.666999999999
$ws.Range("M89").Value = -30006
$ws.Range("N89").Value = -17833.667
$ws.Range("H107").Value = 3119.6572
$ws.Range("I107").Value = 3100.5518
$ws.Range("K107").Value = 3100.5518
$ws.Range("M107").Value = -1180.5518
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H129").Value = 1082.3846
$ws.Range("I129").Value = 1006.4545
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 3019.3635
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = 1980.6365
$ws.Range("N129").Value = -14500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 30001
$ws.Range("I31").Value = 30001
$ws.Range("K31").Value = 30001
$ws.Range("M31").Value = -29707
$ws.Range("H32").Value = 223097.3
$ws.Range("I32").Value = 225832.8
$ws.Range("K32").Value = 225832.8
$ws.Range("M32").Value = -225545.8
$ws.Range("H74").Value = 5445.95
$ws.Range("I74").Value = 2954.3076
$ws.Range("J74").Value = 21641.625
$ws.Range("K74").Value = 2954.3076
$ws.Range("L74").Value = 21641.625
$ws.Range("M74").Value = -2080.3076
$ws.Range("N74").Value = -23389.625
$ws.Range("H77").Value = 5445.95
$ws.Range("I77").Value = 2954.3076
$ws.Range("J77").Value = 21641.625
$ws.Range("K77").Value = 14771.538
$ws.Range("L77").Value = 108208.125
$ws.Range("M77").Value = -10403.538
$ws.Range("N77").Value = -116944.125
$ws.Range("H122").Value = 2544
$ws.Range("I122").Value = 2396.3635
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 7189.0905
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -4739.0905
$ws.Range("N122").Value = -13750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2706.8572
$ws.Range("I86").Value = 2210.6667
$ws.Range("K86").Value = 2210.6667
$ws.Range("M86").Value = -1087.6667
$ws.Range("H89").Value = 2706.8572
$ws.Range("I89").Value = 2210.6667
$ws.Range("K89").Value = 11053.3335
$ws.Range("M89").Value = -5437.333500000001
$ws.Range("H102").Value = 34501.4
$ws.Range("I102").Value = 26876.75
$ws.Range("K102").Value = 26876.75
$ws.Range("M102").Value = -23631.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 51999.5
$ws.Range("I18").Value = 39999
$ws.Range("J18").Value = 64000
$ws.Range("K18").Value = 39999
$ws.Range("L18").Value = 64000
$ws.Range("M18").Value = -39769
$ws.Range("N18").Value = -64460
$ws.Range("H31").Value = 2608.4375
$ws.Range("I31").Value = 2249
$ws.Range("K31").Value = 2249
$ws.Range("M31").Value = -1954
$ws.Range("H34").Value = 2608.4375
$ws.Range("I34").Value = 2249
$ws.Range("K34").Value = 2249
$ws.Range("M34").Value = -2047
$ws.Range("H41").Value = 24500
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 24500
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 24500
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -25356
$ws.Range("H53").Value = 52536.8
$ws.Range("I53").Value = 55000
$ws.Range("J53").Value = 42684
$ws.Range("K53").Value = 55000
$ws.Range("L53").Value = 42684
$ws.Range("M53").Value = -54393
$ws.Range("N53").Value = -43898
$ws.Range("H60").Value = 22599.834
$ws.Range("I60").Value = 17000
$ws.Range("J60").Value = 25399.75
$ws.Range("K60").Value = 17000
$ws.Range("L60").Value = 25399.75
$ws.Range("M60").Value = -16489
$ws.Range("N60").Value = -26421.75
$ws.Range("H99").Value = 14690.6875
$ws.Range("J99").Value = 2581.6667
$ws.Range("L99").Value = 2581.6667
$ws.Range("N99").Value = -5577.6667
$ws.Range("H114").Value = 44666.668
$ws.Range("J114").Value = 44666.668
$ws.Range("L114").Value = 44666.668
$ws.Range("N114").Value = -53344.668
$ws.Range("H122").Value = 9819.546
$ws.Range("I122").Value = 2086.8518
$ws.Range("K122").Value = 6260.555399999999
$ws.Range("M122").Value = -3810.555399999999
$ws.Range("H126").Value = 14690.6875
$ws.Range("J126").Value = 2581.6667
$ws.Range("L126").Value = 7745.000100000001
$ws.Range("N126").Value = -12685.0001
$ws.Range("H132").Value = 2605.4243
$ws.Range("I132").Value = 2673.037
$ws.Range("J132").Value = 2301.1667
$ws.Range("K132").Value = 8019.110999999999
$ws.Range("L132").Value = 6903.500100000001
$ws.Range("M132").Value = -5489.110999999999
$ws.Range("N132").Value = -11963.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 203
$ws.Range("I2").Value = 198.66667
$ws.Range("K2").Value = 1192.00002
$ws.Range("M2").Value = -1079.00002
$ws.Range("H17").Value = 1819.3334
$ws.Range("I17").Value = 3458.3333
$ws.Range("J17").Value = 999.8333
$ws.Range("K17").Value = 10374.9999
$ws.Range("L17").Value = 2999.4999
$ws.Range("M17").Value = -10205.9999
$ws.Range("N17").Value = -3337.4999
$ws.Range("H34").Value = 2400.2144
$ws.Range("J34").Value = 2827.3
$ws.Range("L34").Value = 8481.900000000001
$ws.Range("N34").Value = -8649.900000000001
$ws.Range("H39").Value = 120797.72
$ws.Range("J39").Value = 151181.81
$ws.Range("L39").Value = 453545.43
$ws.Range("N39").Value = -454133.43
$ws.Range("H55").Value = 94447460
$ws.Range("J55").Value = 2004715.6
$ws.Range("L55").Value = 6014146.800000001
$ws.Range("N55").Value = -6014500.800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12501629
$ws.Range("I43").Value = 12501629
$ws.Range("K43").Value = 12501629
$ws.Range("M43").Value = -12501478
$ws.Range("H80").Value = 1999
$ws.Range("I80").Value = 1999
$ws.Range("J80").Value = 1999
$ws.Range("K80").Value = 1999
$ws.Range("L80").Value = 1999
$ws.Range("M80").Value = -1001
$ws.Range("N80").Value = -3995
$ws.Range("H83").Value = 1999
$ws.Range("I83").Value = 1999
$ws.Range("J83").Value = 1999
$ws.Range("K83").Value = 9995
$ws.Range("L83").Value = 9995
$ws.Range("M83").Value = -5003
$ws.Range("N83").Value = -19979
$ws.Range("H102").Value = 2554.9565
$ws.Range("J102").Value = 1181.8
$ws.Range("L102").Value = 1181.8
$ws.Range("N102").Value = -4425.8
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988
$ws.Range("H122").Value = 4354.0557
$ws.Range("I122").Value = 4434.4546
$ws.Range("K122").Value = 13303.3638
$ws.Range("M122").Value = -10853.3638
$ws.Range("H126").Value = 2853.9
$ws.Range("I126").Value = 2846.2
$ws.Range("K126").Value = 8538.599999999999
$ws.Range("M126").Value = -6068.599999999999
$ws.Range("H140").Value = 99390
$ws.Range("J140").Value = 99390
$ws.Range("L140").Value = 99390
$ws.Range("N140").Value = -109750
$ws.Range("H141").Value = 73086.5
$ws.Range("J141").Value = 73086.5
$ws.Range("L141").Value = 73086.5
$ws.Range("N141").Value = -83446.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1804.909
$ws.Range("I16").Value = 1805.8889
$ws.Range("K16").Value = 1805.8889
$ws.Range("M16").Value = -1635.8889
$ws.Range("H22").Value = 2731.92
$ws.Range("I22").Value = 2551.1667
$ws.Range("J22").Value = 2898.7693
$ws.Range("K22").Value = 2551.1667
$ws.Range("L22").Value = 2898.7693
$ws.Range("M22").Value = -2256.1667
$ws.Range("N22").Value = -3488.7693
$ws.Range("H27").Value = 2731.92
$ws.Range("I27").Value = 2551.1667
$ws.Range("J27").Value = 2898.7693
$ws.Range("K27").Value = 2551.1667
$ws.Range("L27").Value = 2898.7693
$ws.Range("M27").Value = -2444.1667
$ws.Range("N27").Value = -3112.7693
$ws.Range("H82").Value = 3122.4666
$ws.Range("J82").Value = 2712.5
$ws.Range("L82").Value = 2712.5
$ws.Range("N82").Value = -3434.5
$ws.Range("H85").Value = 3122.4666
$ws.Range("J85").Value = 2712.5
$ws.Range("L85").Value = 2712.5
$ws.Range("N85").Value = -5208.5
$ws.Range("H132").Value = 3129.56
$ws.Range("J132").Value = 3897.6155
$ws.Range("L132").Value = 11692.8465
$ws.Range("N132").Value = -16752.8465
$ws.Range("H133").Value = 79536.5
$ws.Range("J133").Value = 88777
$ws.Range("L133").Value = 88777
$ws.Range("N133").Value = -93837

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 17018334
$ws.Range("I49").Value = 25500000
$ws.Range("K49").Value = 25500000
$ws.Range("M49").Value = -25499770
$ws.Range("H81").Value = 49379.453
$ws.Range("I81").Value = 4066
$ws.Range("K81").Value = 8132
$ws.Range("M81").Value = -7071
$ws.Range("H84").Value = 49379.453
$ws.Range("I84").Value = 4066
$ws.Range("K84").Value = 40660
$ws.Range("M84").Value = -35356
